$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.125.70'
$ws.Range("E2").Value = '  -3.33%  '
$ws.Range("D3").Value = '1.605.80'
$ws.Range("E3").Value = '  -2.93%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("D6").Value = '''302.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.94%  '
$ws.Range("D7").Value = '''0.3792'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.57%  '
$ws.Range("D8").Value = '''0.3668'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.93%  '
$ws.Range("D9").Value = '''50.31'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.32%  '
$ws.Range("D10").Value = '''1.276'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.28%  '
$ws.Range("D11").Value = '''0.08160'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.82%  '
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").Value = '''22.92'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.49%  '
$ws.Range("D14").Value = '''6.634'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.92%  '
$ws.Range("E15").Value = '  -4.03%  '
$ws.Range("D16").Value = '''7.422'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.02%  '
$ws.Range("D17").Value = '1.602.40'
$ws.Range("E17").Value = '  -3.68%  '
$ws.Range("D18").Value = '''92.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("D19").Value = '''0.06873'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("D20").Value = '''18.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.53%  '
$ws.Range("D21").Value = '''6.618'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.36%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '''1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '''13.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.45%  '
$ws.Range("B24").Value = 'WrappedBTC'
$ws.Range("C24").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D24").Value = '23.132.75'
$ws.Range("E24").Value = '  -3.28%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '''2.360'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.02%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '''2.820'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.84%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''21.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.87%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '''150.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.52%  '
$ws.Range("B29").Value = 'HuobiToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D29").Value = '''5.285'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.78%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '''134.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.45%  '
$ws.Range("B31").Value = 'WEMIXTOKEN'
$ws.Range("C31").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D31").Value = '''2.387'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.15%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''6.896'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -11.62%  '
$ws.Range("B33").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C33").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D33").Value = '1.781.28'
$ws.Range("E33").Value = '  -3.29%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''0.9643'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.31%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.07748'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.16%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '''10.45'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.13%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '''6.326'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.13%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.02742'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.93%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '''0.2559'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.35%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '''0.08903'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.57%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '''1.371'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.84%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '''0.7113'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.94%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '''12.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.66%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '''15.35'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.01%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '''0.6652'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.07%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '''2.332'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.00%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '''0.9992'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Value = '''4.005'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.55%  '
$ws.Range("B49").Value = 'Flow'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D49").Value = '''1.255'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.95%  '
$ws.Range("D50").Value = '''132.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.07954'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.01%  '
